# Adds "xOffsetsUnits" and "yOffsetsUnits" columns (with per-row unit values)
# to the DataCombined sheet, and removes the now-unused "tagPrefix" column
# header on the plotGrids sheet.

$wb = $excel.ActiveWorkbook

# --- Sheet "DataCombined" (1st sheet) -------------------------------------
$ws1 = $wb.Worksheets.Item(1)

# Original header layout:  ... H=xOffsets  I=yOffsets        J=xScaleFactors  K=yScaleFactors
# New header layout:       ... H=xOffsets  I=xOffsetsUnits  J=yOffsets  K=yOffsetsUnits  L=xScaleFactors  M=yScaleFactors

# Insert a column before the existing "yOffsets" column (col I) for the new
# "xOffsetsUnits" column.
$ws1.Columns.Item(9).Insert()

# Insert a second column before the (now shifted) "xScaleFactors" column
# (col K) for the new "yOffsetsUnits" column.
$ws1.Columns.Item(11).Insert()

$ws1.Range("I1").Value = "xOffsetsUnits"
$ws1.Range("K1").Value = "yOffsetsUnits"

# Populate the unit values for the two existing data rows.
$ws1.Range("H2").Value = 1
$ws1.Range("I2").Value = "h"

$ws1.Range("H3").Value = 1
$ws1.Range("I3").Value = "min"

# --- Sheet "plotGrids" (3rd sheet) -----------------------------------------
$ws3 = $wb.Worksheets.Item(3)

# Select D1 before clearing it so the saved view matches the new selection.
$ws3.Range("D1").Select()

# The "tagPrefix" column is no longer used; remove its header value.
$ws3.Range("D1").ClearContents()

# --- Final view/selection state ---------------------------------------------
$ws1.Activate()
$ws1.Range("I4").Select()
